$d = $word.ActiveDocument

# Remove the trailing "Ver no Jupiter ..." and copyright paragraphs, along
# with the blank paragraph that separates them from the bibliography entry,
# leaving the bibliography text and the final blank/page-break paragraphs
# untouched.
$start = $d.Paragraphs.Item(41).Range.Start
$end = $d.Paragraphs.Item(43).Range.End

$r = $d.Range($start, $end)
$r.Delete()
